# Update the Terminal Handling Database transport calculator:
# rework the POD pricing tiers (Rotterdam motorcycle surcharge and the
# whole Varna price column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rotterdam (POD) rows: MOTORCYCLE price 700 -> 600 ---
foreach ($r in @(6, 11, 16, 21, 26, 31, 36)) {
    $ws.Cells.Item($r, 4).Value = 600
}

# --- Varna (POD) rows: CAR / SUV / LARGE SUV / PICKUP price 700 -> 400 ---
foreach ($r in 37..70) {
    if ((($r - 37) % 5) -ne 4) {
        $ws.Cells.Item($r, 4).Value = 400
    }
}

# --- Varna (POD) rows: MOTORCYCLE price 700 -> 300 ---
foreach ($r in @(41, 46, 51, 56, 61, 66, 71)) {
    $ws.Cells.Item($r, 4).Value = 300
}

# Reflect the last on-screen selection left by the editor.
$ws.Range("H67").Select() | Out-Null
